$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The export now only reports the latest (2014) figure, so drop the left
# border that used to visually separate the "Area" row label from the
# multi-year grid.
$ws.Range("A6").Borders.Item(7).LineStyle = -4142

# Remove the columns that held the 1989 and 2002 census figures; only the
# 2014 column remains (shifts into column B).
$ws.Range("B:C").Delete()

# Remove the now-unused spacer row that used to separate the title block
# from the "(sq. km)" row.
$ws.Range("3:3").Delete()

# The "(according to the population census data)" subtitle is removed from
# the export entirely.
$ws.Range("A2").Clear()

# The refreshed export uses taller header/data rows.
$ws.Range("1:5").RowHeight = 20.1

Write-Output "done"
